$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Text = "50÷2="
$cell = $t.Cell(1, 2)
$cell.Range.Text = "45÷9="
$cell = $t.Cell(1, 3)
$cell.Range.Text = "36÷4="
$cell = $t.Cell(1, 4)
$cell.Range.Text = "74÷7="
$cell = $t.Cell(1, 5)
$cell.Range.Text = "63÷5="
$cell = $t.Cell(5, 1)
$cell.Range.Text = "40÷8="
$cell = $t.Cell(5, 2)
$cell.Range.Text = "49÷2="
$cell = $t.Cell(5, 3)
$cell.Range.Text = "18÷7="
$cell = $t.Cell(5, 4)
$cell.Range.Text = "35÷8="
$cell = $t.Cell(5, 5)
$cell.Range.Text = "81÷5="
$cell = $t.Cell(9, 1)
$cell.Range.Text = "60÷3="
$cell = $t.Cell(9, 2)
$cell.Range.Text = "29÷4="
$cell = $t.Cell(9, 3)
$cell.Range.Text = "86÷7="
$cell = $t.Cell(9, 4)
$cell.Range.Text = "56÷2="
$cell = $t.Cell(9, 5)
$cell.Range.Text = "40÷6="
$cell = $t.Cell(13, 1)
$cell.Range.Text = "23÷2="
$cell = $t.Cell(13, 2)
$cell.Range.Text = "23÷4="
$cell = $t.Cell(13, 3)
$cell.Range.Text = "76÷7="
$cell = $t.Cell(13, 4)
$cell.Range.Text = "24÷3="
$cell = $t.Cell(13, 5)
$cell.Range.Text = "57÷5="
$cell = $t.Cell(17, 1)
$cell.Range.Text = "39÷8="
$cell = $t.Cell(17, 2)
$cell.Range.Text = "42÷6="
$cell = $t.Cell(17, 3)
$cell.Range.Text = "68÷9="
$cell = $t.Cell(17, 4)
$cell.Range.Text = "48÷3="
$cell = $t.Cell(17, 5)
$cell.Range.Text = "99÷2="
